$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 191 (shifts existing rows 191-296 down to 193-298,
# carrying along formatting such as the date style on column D).
$ws.Rows("191:192").Insert()

# New row 191 ("Primera" quality, week of 2022-04-29 = serial 44680)
$ws.Range("A191").Value = 1
$ws.Range("B191").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C191").Value = "Arica y Parinacota"
$ws.Range("D191").Value = 44680
$ws.Range("E191").Value = 15
$ws.Range("F191").Value = 100112043
$ws.Range("G191").Value = "Pepino ensalada"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 130
$ws.Range("K191").Value = 12000
$ws.Range("L191").Value = 13000
$ws.Range("M191").Value = 12500
$ws.Range("N191").Value = "$/caja 70 unidades"
$ws.Range("O191").Value = "Región de Arica y Parinacota"
$ws.Range("P191").Value = 179
$ws.Range("Q191").Value = 70
$ws.Range("R191").Value = "Hortaliza"

# New row 192 ("Segunda" quality, same week)
$ws.Range("A192").Value = 1
$ws.Range("B192").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C192").Value = "Arica y Parinacota"
$ws.Range("D192").Value = 44680
$ws.Range("E192").Value = 15
$ws.Range("F192").Value = 100112043
$ws.Range("G192").Value = "Pepino ensalada"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Segunda"
$ws.Range("J192").Value = 160
$ws.Range("K192").Value = 10000
$ws.Range("L192").Value = 11000
$ws.Range("M192").Value = 10500
$ws.Range("N192").Value = "$/caja 100 unidades"
$ws.Range("O192").Value = "Región de Arica y Parinacota"
$ws.Range("P192").Value = 105
$ws.Range("Q192").Value = 100
$ws.Range("R192").Value = "Hortaliza"
